# Swap the "Office Theme" and "Integral" theme colour schemes.
#
# The deck ships two theme parts:
#   - the presentation's active theme (used by the slide master / all
#     slides) currently holds the "Integral" / Red Violet palette.
#   - the notes-master theme currently holds the plain "Office Theme"
#     palette.
#
# The target edit swaps their contents so the active theme becomes the
# Office palette and the (otherwise identical) font/format schemes stay
# untouched — only the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) differ between the two themes.
#
# PowerPoint's object model exposes all 12 slots through
# Slide.ThemeColorScheme(1..12).RGB, which writes straight back into the
# active theme's <a:clrScheme>, so we drive the swap through that.

$p = $ppt.ActivePresentation

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as 0xRRGGBB, in theme-colour-slot order.
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $red = ($hex -shr 16) -band 0xFF
    $green = ($hex -shr 8) -band 0xFF
    $blue = $hex -band 0xFF
    # PowerPoint RGB() packing is 0x00BBGGRR.
    $comRgb = ($blue * 65536) + ($green * 256) + $red
    $themeColors.Item($i).RGB = $comRgb
}
